$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header section updates ---
$ws.Range("E11").Value = 819336
$ws.Range("C13").Value = 8
$ws.Range("F13").Value = 21

# --- Insert 19 new rows after row 18 (row 19 becomes the insertion point), shifting
#     the blank gap + footer rows down so the footer lands on rows 42/43 ---
for ($i = 0; $i -lt 19; $i++) {
    $ws.Rows.Item(19).Insert()
}

# --- Re-apply the "normal" data-row formatting (copied from row 17) to the newly
#     inserted rows 19-36, and the "last row" formatting (copied from row 18) to the
#     new final data row 37 ---
$ws.Range("B17:J17").Copy()
$ws.Range("B19:J36").PasteSpecial(-4122)
$ws.Range("B18:J18").Copy()
$ws.Range("B37:J37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the worker / period data rows ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "9097278"
$ws.Range("D16").Value = "ALVARO CASTRO PERIÑAN"
$ws.Range("E16").Value = "2506"
$ws.Range("F16").Value = 31249
$ws.Range("G16").Value = 781242

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73151464"
$ws.Range("D17").Value = "JAIRO VARGAS TRUJILLO"
$ws.Range("E17").Value = "1806"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 781242

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73151464"
$ws.Range("D18").Value = "JAIRO VARGAS TRUJILLO"
$ws.Range("E18").Value = "1805"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 781242

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73151464"
$ws.Range("D19").Value = "JAIRO VARGAS TRUJILLO"
$ws.Range("E19").Value = "1804"
$ws.Range("F19").Value = 31249
$ws.Range("G19").Value = 781242

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "73569467"
$ws.Range("D20").Value = "OSCAR MARRIAGA URUETA"
$ws.Range("E20").Value = "2211"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 877803

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "73569467"
$ws.Range("D21").Value = "OSCAR MARRIAGA URUETA"
$ws.Range("E21").Value = "2210"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 877803

$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "73569467"
$ws.Range("D22").Value = "OSCAR MARRIAGA URUETA"
$ws.Range("E22").Value = "2209"
$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 877803

$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "5110438"
$ws.Range("D23").Value = "YESID DAVILA CONTRERAS"
$ws.Range("E23").Value = "2206"
$ws.Range("F23").Value = 40000
$ws.Range("G23").Value = 1000000

$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "5110438"
$ws.Range("D24").Value = "YESID DAVILA CONTRERAS"
$ws.Range("E24").Value = "2107"
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = 1000000

$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "5110438"
$ws.Range("D25").Value = "YESID DAVILA CONTRERAS"
$ws.Range("E25").Value = "2106"
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = 1000000

$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "5110438"
$ws.Range("D26").Value = "YESID DAVILA CONTRERAS"
$ws.Range("E26").Value = "2105"
$ws.Range("F26").Value = 36341
$ws.Range("G26").Value = 1000000

$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "5110438"
$ws.Range("D27").Value = "YESID DAVILA CONTRERAS"
$ws.Range("E27").Value = "2104"
$ws.Range("F27").Value = 36341
$ws.Range("G27").Value = 1000000

$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "5110438"
$ws.Range("D28").Value = "YESID DAVILA CONTRERAS"
$ws.Range("E28").Value = "2103"
$ws.Range("F28").Value = 36341
$ws.Range("G28").Value = 1000000

$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "5110438"
$ws.Range("D29").Value = "YESID DAVILA CONTRERAS"
$ws.Range("E29").Value = "2102"
$ws.Range("F29").Value = 36341
$ws.Range("G29").Value = 1000000

$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "73185723"
$ws.Range("D30").Value = "DIXON CORREA GONZALEZ"
$ws.Range("E30").Value = "2204"
$ws.Range("F30").Value = 36341
$ws.Range("G30").Value = 908526

$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "73185723"
$ws.Range("D31").Value = "DIXON CORREA GONZALEZ"
$ws.Range("E31").Value = "2203"
$ws.Range("F31").Value = 36341
$ws.Range("G31").Value = 908526

$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "45591343"
$ws.Range("D32").Value = "YINCETH GONZALEZ GONZALEZ"
$ws.Range("E32").Value = "2208"
$ws.Range("F32").Value = 36341
$ws.Range("G32").Value = 908526

$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "45591343"
$ws.Range("D33").Value = "YINCETH GONZALEZ GONZALEZ"
$ws.Range("E33").Value = "2207"
$ws.Range("F33").Value = 36341
$ws.Range("G33").Value = 908526

$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "45591343"
$ws.Range("D34").Value = "YINCETH GONZALEZ GONZALEZ"
$ws.Range("E34").Value = "2206"
$ws.Range("F34").Value = 36341
$ws.Range("G34").Value = 908526

$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "1143333093"
$ws.Range("D35").Value = "RAMIRO MIGUEL CASTRO BELLO"
$ws.Range("E35").Value = "1810"
$ws.Range("F35").Value = 31249
$ws.Range("G35").Value = 781242

$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "1143333093"
$ws.Range("D36").Value = "RAMIRO MIGUEL CASTRO BELLO"
$ws.Range("E36").Value = "1809"
$ws.Range("F36").Value = 31249
$ws.Range("G36").Value = 781242

$ws.Range("B37").Value = "CC"
$ws.Range("C37").Value = "1007314509"
$ws.Range("D37").Value = "JAYBER MORALES MALDONADO"
$ws.Range("E37").Value = "2305"
$ws.Range("F37").Value = 46400
$ws.Range("G37").Value = 1160000
